$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "final changes of browser size" - refresh the Job# list (column B, rows 2-11)
# with the latest batch of job numbers. Values are numeric-looking text, so
# they must stay stored as text (not auto-converted to numbers) - write them
# with a leading quote (forces text) and then strip the resulting format
# back off so the cell keeps its original (default) style.
$newValues = @(
    "32297175",
    "32297176",
    "32297178",
    "32297180",
    "32297186",
    "32297193",
    "32297201",
    "32297209",
    "32297242",
    "32297272"
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = "'" + $newValues[$i]
    $cell.ClearFormats()
}
